$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final sorted language/value data (descending by value), rows 2-21
$data = @(
    @("English", 26.61007214192699),
    @("Spanish", 8.150911593408372),
    @("Japanese", 8.102405126244054),
    @("German", 6.515693814052288),
    @("Chinese", 6.082060391213726),
    @("Arabic", 4.810120066892289),
    @("Russian", 4.449994812542799),
    @("Portuguese", 3.812147723581421),
    @("French", 3.751023817510207),
    @("Italian", 3.640371497933826),
    @("Malay-Indonesian", 2.58744147126189),
    @("Dutch", 1.662381400968103),
    @("Persian", 1.514225090043548),
    @("Turkish", 1.431159564225959),
    @("Korean", 1.267282551984827),
    @("Thai", 0.9640984205142802),
    @("Polish", 0.7839479923158847),
    @("Urdu", 0.774848351799184),
    @("Swedish", 0.514951564349399),
    @("Bengali", 0.3971542565533244)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the two rows that are no longer part of the data (previously rows 22-23)
$ws.Range("A22:B23").Clear()
